# Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy the header style (bold, centered, bordered) from the "Weekly Quantity" sheet
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$data = New-Object 'object[,]' 31,4
$data[0,0] = 44934.99999999999
$data[0,1] = 0
$data[0,2] = -182.1789849798427
$data[0,3] = 67.90769153989369
$data[1,0] = 44997.99999999999
$data[1,1] = 0
$data[1,2] = -133.9517447626357
$data[1,3] = 106.8505465216911
$data[2,0] = 45011.99999999999
$data[2,1] = 0
$data[2,2] = -127.1180501492676
$data[2,3] = 118.0634441594494
$data[3,0] = 45018.99999999999
$data[3,1] = 0
$data[3,2] = -132.2442309535334
$data[3,3] = 119.1821479661968
$data[4,0] = 45032.99999999999
$data[4,1] = 11
$data[4,2] = -110.1307669824207
$data[4,3] = 136.9115789166458
$data[5,0] = 45060.99999999999
$data[5,1] = 33
$data[5,2] = -90.72937200834735
$data[5,3] = 155.2440568966181
$data[6,0] = 45095.99999999999
$data[6,1] = 60
$data[6,2] = -65.77322278821403
$data[6,3] = 178.2774383335845
$data[7,0] = 45116.99999999999
$data[7,1] = 77
$data[7,2] = -58.73770244564246
$data[7,3] = 196.7360556944332
$data[8,0] = 45123.99999999999
$data[8,1] = 83
$data[8,2] = -33.53912286464817
$data[8,3] = 205.3916609430989
$data[9,0] = 45130.99999999999
$data[9,1] = 88
$data[9,2] = -30.83210786859244
$data[9,3] = 207.2444872243296
$data[10,0] = 45137.99999999999
$data[10,1] = 94
$data[10,2] = -29.10717276702929
$data[10,3] = 212.5379526071794
$data[11,0] = 45144.99999999999
$data[11,1] = 99
$data[11,2] = -22.61116942871281
$data[11,3] = 225.6861498969523
$data[12,0] = 45151.99999999999
$data[12,1] = 105
$data[12,2] = -19.55960910367768
$data[12,3] = 224.4577732309023
$data[13,0] = 45158.99999999999
$data[13,1] = 110
$data[13,2] = -5.91737412937778
$data[13,3] = 240.7553799105989
$data[14,0] = 45193.99999999999
$data[14,1] = 138
$data[14,2] = 18.42667073189259
$data[14,3] = 253.8558768002903
$data[15,0] = 45200.99999999999
$data[15,1] = 143
$data[15,2] = 14.33851656608248
$data[15,3] = 269.4110970258065
$data[16,0] = 45214.99999999999
$data[16,1] = 155
$data[16,2] = 35.76782734346265
$data[16,3] = 292.5703900654796
$data[17,0] = 45221.99999999999
$data[17,1] = 160
$data[17,2] = 38.86011423857966
$data[17,3] = 283.8874390391393
$data[18,0] = 45228.99999999999
$data[18,1] = 166
$data[18,2] = 37.43826982815676
$data[18,3] = 285.3969480753717
$data[19,0] = 45235.99999999999
$data[19,1] = 171
$data[19,2] = 52.64274233264238
$data[19,3] = 299.8404436364607
$data[20,0] = 45242.99999999999
$data[20,1] = 177
$data[20,2] = 50.91553364620498
$data[20,3] = 290.9926192895553
$data[21,0] = 45249.99999999999
$data[21,1] = 182
$data[21,2] = 53.2470242761696
$data[21,3] = 308.7127161013501
$data[22,0] = 45263.99999999999
$data[22,1] = 193
$data[22,2] = 66.04885947992175
$data[22,3] = 310.0386353567549
$data[23,0] = 45270.99999999999
$data[23,1] = 199
$data[23,2] = 68.68683428151289
$data[23,3] = 315.3108809941215
$data[24,0] = 45277.99999999999
$data[24,1] = 204
$data[24,2] = 79.131924147993
$data[24,3] = 326.1384203178901
$data[25,0] = 45284.99999999999
$data[25,1] = 210
$data[25,2] = 91.87796413633042
$data[25,3] = 343.8092656045754
$data[26,0] = 45291.99999999999
$data[26,1] = 215
$data[26,2] = 94.66921847006616
$data[26,3] = 335.1815206533572
$data[27,0] = 45298.99999999999
$data[27,1] = 221
$data[27,2] = 105.7249632288048
$data[27,3] = 352.2585876650641
$data[28,0] = 45305.99999999999
$data[28,1] = 226
$data[28,2] = 110.4040440094482
$data[28,3] = 347.6673332392307
$data[29,0] = 45312.99999999999
$data[29,1] = 232
$data[29,2] = 115.1150094188369
$data[29,3] = 353.4800586169916
$data[30,0] = 45319.99999999999
$data[30,1] = 237
$data[30,2] = 115.9026025270062
$data[30,3] = 363.3745403301159
$newSheet.Range("A2:D32").Value = $data

# Copy the date-formatted style used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A32").PasteSpecial(-4122)

$excel.CutCopyMode = 0
